# Update cryptos list - automated GitHub Actions style refresh of price/volume data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value while forcing text storage so that Excel does not
# auto-convert numeric-looking strings (e.g. '64.80', '1.00', '0.0480') into
# actual numbers, which would silently drop significant trailing/leading zeros.
# The cell's original style is preserved (format is only toggled transiently).
function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.913.01"
$ws.Range("E2").Value = "  -0.15%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.630.17"
$ws.Range("E3").Value = "  -0.07%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.03%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "211.84"
$ws.Range("E5").Value = "  +0.01%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -1.34%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.04%  "

# Row 8 - Solana
Set-TextValue $ws.Range("D8") "23.23"
$ws.Range("E8").Value = "  -0.77%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -0.08%  "

# Row 10 - Dogecoin
Set-TextValue $ws.Range("D10") "0.0607"
$ws.Range("E10").Value = "  -1.05%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.21%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.861.47"
$ws.Range("E12").Value = "  -0.08%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.625.93"
$ws.Range("E13").Value = "  -0.24%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -1.00%  "

# Row 15 - Polygon
Set-TextValue $ws.Range("D15") "0.555"
$ws.Range("E15").Value = "  -1.26%  "

# Row 16 - Litecoin
Set-TextValue $ws.Range("D16") "64.80"
$ws.Range("E16").Value = "  -1.20%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "27.921.40"
$ws.Range("E17").Value = "  -0.09%  "

# Row 18 - BitcoinCash
Set-TextValue $ws.Range("D18") "227.85"
$ws.Range("E18").Value = "  -1.29%  "

# Row 19 - Chainlink
Set-TextValue $ws.Range("D19") "7.61"
$ws.Range("E19").Value = "  -0.61%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  -0.94%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.08%  "

# Row 22 - Uniswap
Set-TextValue $ws.Range("D22") "4.34"
$ws.Range("E22").Value = "  -0.09%  "

# Row 23 - Avalanche
Set-TextValue $ws.Range("D23") "9.99"
$ws.Range("E23").Value = "  -3.58%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +1.64%  "

# Row 25 - Monero
Set-TextValue $ws.Range("D25") "154.72"
$ws.Range("E25").Value = "  -0.17%  "

# Row 26 - Cosmos
$ws.Range("E26").Value = "  -0.24%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  -0.39%  "

# Row 28 - BinanceUSD
$ws.Range("E28").Value = "  +0.05%  "

# Row 29 - EthereumClassic
Set-TextValue $ws.Range("D29") "15.39"
$ws.Range("E29").Value = "  -1.25%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -0.36%  "

# Row 31 - Hedera
Set-TextValue $ws.Range("D31") "0.0480"
$ws.Range("E31").Value = "  -0.33%  "

# Row 32 - Filecoin
Set-TextValue $ws.Range("D32") "3.39"
$ws.Range("E32").Value = "  -0.17%  "

# Row 33 - Maker
$ws.Range("D33").Value = "1.417.02"
$ws.Range("E33").Value = "  +1.07%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  +0.95%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  +2.60%  "

# Row 36 - TrustWalletToken
Set-TextValue $ws.Range("D36") "1.00"
$ws.Range("E36").Value = "  -1.68%  "

# Row 37 - HuobiToken
$ws.Range("E37").Value = "  -1.12%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  -0.98%  "

# Row 39 - ImmutableX
Set-TextValue $ws.Range("D39") "0.554"
$ws.Range("E39").Value = "  -0.54%  "

# Row 41 - WEMIXToken
$ws.Range("E41").Value = "  -2.22%  "

# Row 42 - RenderToken -> Aave (swap)
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D42") "65.77"
$ws.Range("E42").Value = "  -1.23%  "

# Row 43 - Aave -> RenderToken (swap)
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D43") "1.82"
$ws.Range("E43").Value = "  -1.27%  "

# Row 44 - FraxShare
$ws.Range("E44").Value = "  -0.93%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").Value = "1.770.36"
$ws.Range("E45").Value = "  -0.15%  "

# Row 46 - MXToken
Set-TextValue $ws.Range("D46") "2.12"
$ws.Range("E46").Value = "  -3.57%  "

# Row 47 - Quant
Set-TextValue $ws.Range("D47") "88.58"
$ws.Range("E47").Value = "  +0.43%  "

# Row 48 - BabyDogeCoin -> Algorand (swap)
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D48") "0.101"
$ws.Range("E48").Value = "  +0.42%  "

# Row 49 - Algorand -> BabyDogeCoin (swap)
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0102"
$ws.Range("E49").Value = "  -2.92%  "

# Row 50 - Cronos
$ws.Range("E50").Value = "  -0.36%  "

# Row 51 - EnergySwap
Set-TextValue $ws.Range("D51") "7.61"
$ws.Range("E51").Value = "  +0.79%  "

